$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.546376466751099
$ws.Range("B1").Value = 2.274316072463989
$ws.Range("C1").Value = 4.341536521911621
$ws.Range("D1").Value = 1.759642958641052
$ws.Range("E1").Value = 0.8167190551757812
